$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("02022019")

# Fill in real coordinates that had been left as 0/0 placeholders.
$ws.Range("F4").Value = 37.870640000000002
$ws.Range("G4").Value = -122.25675099999999

$ws.Range("F10").Value = 37.875394
$ws.Range("G10").Value = -122.255144

# F10/G10 previously carried the "missing data" red-highlight style (same
# style as F4/G4). Now that real coordinates are filled in, clear that
# highlight by reusing the plain style already used by other filled-in
# coordinate cells (e.g. F3/G3) - copy its formatting only.
$ws.Range("F3").Copy()
$ws.Range("F10:G10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the view: scroll back to the top-left and move the selection to
# the newly-edited F10:G10 range.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F10:G10").Select()
